# Fruta / hortaliza, semanal
# Updates the weekly price records for rows 3-11 (Poroto granado, Agrícola del
# Norte S.A. de Arica): the "Fecha" (D), "Volumen" (J), "Precio mínimo" (K),
# "Precio máximo" (L), "Precio promedio ponderado" (M) and "Precio $/Kg" (P)
# columns are refreshed with the latest weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3  = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
    4  = @{ D = 45062; J = 1700; K = 2800; L = 3000; M = 2900; P = 2900 }
    5  = @{ D = 44883; J = 290;  K = 1400; L = 1500; M = 1434; P = 1434 }
    6  = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    7  = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
    8  = @{ D = 44907; J = 2300; K = 900;  L = 1000; M = 952;  P = 952  }
    9  = @{ D = 44895; J = 200;  K = 1200; L = 1300; M = 1255; P = 1255 }
    10 = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    11 = @{ D = 44893; J = 3300; K = 1200; L = 1300; M = 1261; P = 1261 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J  # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K  # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $vals.L  # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $vals.M  # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals.P  # P: Precio $/Kg
}
